# Refresh the "想去人数" (want-to-go count) column F figures across all
# sheets to the newly scraped values (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 3623
$ws.Range("F5").Value  = 3623
$ws.Range("F6").Value  = 268
$ws.Range("F7").Value  = 5148
$ws.Range("F8").Value  = 541
$ws.Range("F9").Value  = 374
$ws.Range("F10").Value = 202
$ws.Range("F11").Value = 700
$ws.Range("F13").Value = 99
$ws.Range("F19").Value = 158
$ws.Range("F21").Value = 362
$ws.Range("F22").Value = 4938
$ws.Range("F26").Value = 6064
$ws.Range("F29").Value = 3228
$ws.Range("F30").Value = 348
$ws.Range("F31").Value = 717
$ws.Range("F32").Value = 4446
$ws.Range("F35").Value = 142
$ws.Range("F36").Value = 1045
$ws.Range("F40").Value = 878
$ws.Range("F41").Value = 1028
$ws.Range("F42").Value = 2035

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 14

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1125
$ws.Range("F4").Value = 52

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 1125
$ws.Range("F5").Value  = 52
$ws.Range("F7").Value  = 3623
$ws.Range("F8").Value  = 3623
$ws.Range("F9").Value  = 268
$ws.Range("F10").Value = 5148
$ws.Range("F11").Value = 541
$ws.Range("F12").Value = 374
$ws.Range("F13").Value = 202
$ws.Range("F14").Value = 700
$ws.Range("F16").Value = 99
$ws.Range("F23").Value = 158
$ws.Range("F25").Value = 362
$ws.Range("F26").Value = 4938
$ws.Range("F30").Value = 6064
$ws.Range("F33").Value = 3228
$ws.Range("F34").Value = 348
$ws.Range("F35").Value = 717
$ws.Range("F36").Value = 4446
$ws.Range("F40").Value = 142
$ws.Range("F41").Value = 1045
$ws.Range("F45").Value = 878
$ws.Range("F46").Value = 1028
$ws.Range("F47").Value = 14
$ws.Range("F48").Value = 2035
